$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("K3").Value = '[''Belgium'', ''France'']'
$ws.Range("K4").Value = '[''Belgium'', ''France'']'
$ws.Range("K5").Value = '[''Belgium'', ''France'']'
$ws.Range("K6").Value = '[''Belgium'', ''France'']'
$ws.Range("K7").Value = '[''Belgium'', ''France'']'
$ws.Range("K13").Value = '[''Spain'', ''West Germany'']'
$ws.Range("K14").Value = '[''Spain'', ''West Germany'']'
$ws.Range("K15").Value = '[''Spain'', ''West Germany'']'
$ws.Range("K17").Value = '[''Spain'', ''Portugal'']'
$ws.Range("K18").Value = '[''Italy'', ''West Germany'']'
$ws.Range("K19").Value = '[''Italy'', ''West Germany'']'
$ws.Range("K20").Value = '[''Italy'', ''West Germany'']'
$ws.Range("K21").Value = '[''Italy'', ''West Germany'']'
$ws.Range("K22").Value = '[''Italy'', ''West Germany'']'
$ws.Range("K23").Value = '[''Italy'', ''West Germany'']'
$ws.Range("K24").Value = '[''Italy'', ''West Germany'']'
$ws.Range("K25").Value = '[''Soviet Union'', ''Republic of Ireland'']'
$ws.Range("K26").Value = '[''Soviet Union'', ''Republic of Ireland'']'
$ws.Range("K27").Value = '[''Soviet Union'', ''Republic of Ireland'']'
$ws.Range("K28").Value = '[''Soviet Union'', ''Republic of Ireland'']'
$ws.Range("K29").Value = '[''Soviet Union'', ''Republic of Ireland'']'
$ws.Range("K30").Value = '[''Soviet Union'', ''Republic of Ireland'']'
$ws.Range("K31").Value = '[''Soviet Union'', ''Netherlands'']'
$ws.Range("K33").Value = '[''England'', ''Sweden'']'
$ws.Range("K34").Value = '[''England'', ''Sweden'']'
$ws.Range("K49").Value = '[''Bulgaria'', ''Spain'']'
$ws.Range("K56").Value = '[''England'', ''Netherlands'']'
$ws.Range("K57").Value = '[''England'', ''Netherlands'']'
$ws.Range("K58").Value = '[''England'', ''Netherlands'']'
$ws.Range("K59").Value = '[''England'', ''Netherlands'']'
$ws.Range("K60").Value = '[''England'', ''Netherlands'']'
$ws.Range("K61").Value = '[''England'', ''Netherlands'']'
$ws.Range("K62").Value = '[''England'', ''Scotland'']'
$ws.Range("K63").Value = '[''England'', ''Netherlands'']'
$ws.Range("K64").Value = '[''Croatia'', ''Portugal'']'
$ws.Range("K65").Value = '[''Croatia'', ''Portugal'']'
$ws.Range("K66").Value = '[''Croatia'', ''Portugal'']'
$ws.Range("K67").Value = '[''Croatia'', ''Portugal'']'
$ws.Range("K68").Value = '[''Croatia'', ''Portugal'']'
$ws.Range("K69").Value = '[''Croatia'', ''Portugal'']'
$ws.Range("K70").Value = '[''Croatia'', ''Portugal'']'
$ws.Range("K71").Value = '[''Croatia'', ''Portugal'']'
$ws.Range("K72").Value = '[''Croatia'', ''Portugal'']'
$ws.Range("K73").Value = '[''Germany'', ''Czech Republic'']'
$ws.Range("K74").Value = '[''Germany'', ''Czech Republic'']'
$ws.Range("K75").Value = '[''Germany'', ''Czech Republic'']'
$ws.Range("K76").Value = '[''Germany'', ''Czech Republic'']'
$ws.Range("K77").Value = '[''Germany'', ''Czech Republic'']'
$ws.Range("K78").Value = '[''Germany'', ''Czech Republic'']'
$ws.Range("K79").Value = '[''Germany'', ''Italy'']'
$ws.Range("K80").Value = '[''Germany'', ''Czech Republic'']'
$ws.Range("K81").Value = '[''Belgium'', ''Italy'']'
$ws.Range("K82").Value = '[''Belgium'', ''Italy'']'
$ws.Range("K83").Value = '[''Belgium'', ''Italy'']'
$ws.Range("K84").Value = '[''Italy'', ''Turkey'']'
$ws.Range("K85").Value = '[''Italy'', ''Turkey'']'
$ws.Range("K86").Value = '[''Italy'', ''Turkey'']'
$ws.Range("K87").Value = '[''Italy'', ''Turkey'']'
$ws.Range("K88").Value = '[''England'', ''Portugal'']'
$ws.Range("K89").Value = '[''England'', ''Portugal'']'
$ws.Range("K92").Value = '[''England'', ''Portugal'']'
$ws.Range("K93").Value = '[''England'', ''Portugal'']'
$ws.Range("K94").Value = '[''England'', ''Portugal'']'
$ws.Range("K95").Value = '[''England'', ''Portugal'']'
$ws.Range("K96").Value = '[''England'', ''Portugal'']'
$ws.Range("K98").Value = '[''FR Yugoslavia'', ''Norway'']'
$ws.Range("K99").Value = '[''FR Yugoslavia'', ''Norway'']'
$ws.Range("K100").Value = '[''FR Yugoslavia'', ''Norway'']'
$ws.Range("K101").Value = '[''FR Yugoslavia'', ''Norway'']'
$ws.Range("K102").Value = '[''FR Yugoslavia'', ''Norway'']'
$ws.Range("K103").Value = '[''FR Yugoslavia'', ''Norway'']'
$ws.Range("K104").Value = '[''FR Yugoslavia'', ''Norway'']'
$ws.Range("K105").Value = '[''FR Yugoslavia'', ''Norway'']'
$ws.Range("K116").Value = '[''Spain'', ''Greece'']'
$ws.Range("K117").Value = '[''Spain'', ''Greece'']'
$ws.Range("K118").Value = '[''Spain'', ''Greece'']'
$ws.Range("K119").Value = '[''Spain'', ''Greece'']'
$ws.Range("K120").Value = '[''Spain'', ''Greece'']'
$ws.Range("K121").Value = '[''Greece'', ''Portugal'']'
$ws.Range("K122").Value = '[''England'', ''France'']'
$ws.Range("K126").Value = '[''England'', ''France'']'
$ws.Range("K127").Value = '[''England'', ''France'']'
$ws.Range("K128").Value = '[''England'', ''France'']'
$ws.Range("K129").Value = '[''England'', ''France'']'
$ws.Range("K130").Value = '[''England'', ''France'']'
$ws.Range("K131").Value = '[''England'', ''France'']'
$ws.Range("K132").Value = '[''England'', ''France'']'
$ws.Range("K145").Value = '[''Netherlands'', ''Czech Republic'']'
$ws.Range("K146").Value = '[''Netherlands'', ''Czech Republic'']'
$ws.Range("K147").Value = '[''Netherlands'', ''Czech Republic'']'
$ws.Range("K148").Value = '[''Netherlands'', ''Czech Republic'']'
$ws.Range("K149").Value = '[''Czech Republic'', ''Portugal'']'
$ws.Range("K150").Value = '[''Czech Republic'', ''Portugal'']'
$ws.Range("K151").Value = '[''Czech Republic'', ''Portugal'']'
$ws.Range("K152").Value = '[''Czech Republic'', ''Portugal'']'
$ws.Range("K153").Value = '[''Czech Republic'', ''Portugal'']'
$ws.Range("K154").Value = '[''Czech Republic'', ''Portugal'']'
$ws.Range("K155").Value = '[''Czech Republic'', ''Portugal'']'
$ws.Range("K156").Value = '[''Czech Republic'', ''Portugal'']'
$ws.Range("K163").Value = '[''Netherlands'', ''Romania'']'
$ws.Range("K164").Value = '[''Netherlands'', ''Romania'']'
$ws.Range("K165").Value = '[''Netherlands'', ''Romania'']'
$ws.Range("K166").Value = '[''Netherlands'', ''Italy'']'
$ws.Range("K167").Value = '[''Netherlands'', ''Italy'']'
$ws.Range("K168").Value = '[''Netherlands'', ''Italy'']'
$ws.Range("K169").Value = '[''Netherlands'', ''Italy'']'
$ws.Range("K172").Value = '[''Spain'', ''Russia'']'
$ws.Range("K173").Value = '[''Spain'', ''Russia'']'
$ws.Range("K174").Value = '[''Spain'', ''Russia'']'
$ws.Range("K175").Value = '[''Spain'', ''Russia'']'
$ws.Range("K176").Value = '[''Spain'', ''Russia'']'
$ws.Range("K177").Value = '[''Czech Republic'', ''Russia'']'
$ws.Range("K178").Value = '[''Czech Republic'', ''Russia'']'
$ws.Range("K179").Value = '[''Czech Republic'', ''Russia'']'
$ws.Range("K180").Value = '[''Greece'', ''Russia'']'
$ws.Range("K181").Value = '[''Czech Republic'', ''Greece'']'
$ws.Range("K182").Value = '[''Germany'', ''Portugal'']'
$ws.Range("K183").Value = '[''Germany'', ''Denmark'']'
$ws.Range("K184").Value = '[''Germany'', ''Portugal'']'
$ws.Range("K185").Value = '[''Germany'', ''Denmark'']'
$ws.Range("K186").Value = '[''Germany'', ''Portugal'']'
$ws.Range("K187").Value = '[''Germany'', ''Portugal'']'
$ws.Range("K188").Value = '[''Germany'', ''Portugal'']'
$ws.Range("K189").Value = '[''Croatia'', ''Spain'']'
$ws.Range("K190").Value = '[''Croatia'', ''Spain'']'
$ws.Range("K191").Value = '[''Croatia'', ''Spain'']'
$ws.Range("K192").Value = '[''Spain'', ''Italy'']'
$ws.Range("K193").Value = '[''Spain'', ''Italy'']'
$ws.Range("K194").Value = '[''Spain'', ''Italy'']'
$ws.Range("K195").Value = '[''England'', ''France'']'
$ws.Range("K196").Value = '[''England'', ''France'']'
$ws.Range("K197").Value = '[''England'', ''France'']'
$ws.Range("K198").Value = '[''England'', ''France'']'
$ws.Range("K199").Value = '[''England'', ''France'']'
$ws.Range("K200").Value = '[''England'', ''France'']'
$ws.Range("K201").Value = '[''France'', ''Switzerland'', ''Romania'']'
$ws.Range("K202").Value = '[''France'', ''Switzerland'', ''Romania'']'
$ws.Range("K203").Value = '[''France'', ''Switzerland'', ''Romania'']'
$ws.Range("K204").Value = '[''France'', ''Albania'', ''Switzerland'']'
$ws.Range("K205").Value = '[''England'', ''Wales'', ''Slovakia'']'
$ws.Range("K206").Value = '[''England'', ''Wales'', ''Slovakia'']'
$ws.Range("K207").Value = '[''England'', ''Wales'', ''Slovakia'']'
$ws.Range("K208").Value = '[''England'', ''Wales'', ''Slovakia'']'
$ws.Range("K209").Value = '[''England'', ''Wales'', ''Slovakia'']'
$ws.Range("K210").Value = '[''England'', ''Wales'', ''Slovakia'']'
$ws.Range("K211").Value = '[''Poland'', ''Germany'', ''Northern Ireland'']'
$ws.Range("K212").Value = '[''Poland'', ''Germany'', ''Northern Ireland'']'
$ws.Range("K213").Value = '[''Poland'', ''Germany'', ''Northern Ireland'']'
$ws.Range("K214").Value = '[''Poland'', ''Germany'', ''Northern Ireland'']'
$ws.Range("K215").Value = '[''Poland'', ''Germany'', ''Northern Ireland'']'
$ws.Range("K216").Value = '[''Croatia'', ''Spain'', ''Czech Republic'']'
$ws.Range("K217").Value = '[''Croatia'', ''Spain'', ''Czech Republic'']'
$ws.Range("K218").Value = '[''Croatia'', ''Spain'', ''Czech Republic'']'
$ws.Range("K219").Value = '[''Croatia'', ''Spain'', ''Turkey'']'
$ws.Range("K220").Value = '[''Croatia'', ''Spain'', ''Turkey'']'
$ws.Range("K221").Value = '[''Croatia'', ''Spain'', ''Turkey'']'
$ws.Range("K222").Value = '[''Croatia'', ''Spain'', ''Turkey'']'
$ws.Range("K223").Value = '[''Portugal'', ''Hungary'', ''Iceland'']'
$ws.Range("K224").Value = '[''Portugal'', ''Hungary'', ''Iceland'']'
$ws.Range("K225").Value = '[''Portugal'', ''Hungary'', ''Iceland'']'
$ws.Range("K226").Value = '[''Portugal'', ''Hungary'', ''Iceland'']'
$ws.Range("K227").Value = '[''Portugal'', ''Hungary'', ''Iceland'']'
$ws.Range("K228").Value = '[''Portugal'', ''Hungary'', ''Iceland'']'
$ws.Range("K229").Value = '[''Portugal'', ''Hungary'', ''Iceland'']'
$ws.Range("K230").Value = '[''Portugal'', ''Hungary'', ''Iceland'']'
$ws.Range("K231").Value = '[''Iceland'', ''Hungary'', ''Portugal'']'
$ws.Range("K232").Value = '[''Portugal'', ''Hungary'', ''Iceland'']'
$ws.Range("K233").Value = '[''Belgium'', ''Italy'', ''Sweden'']'
$ws.Range("K234").Value = '[''Belgium'', ''Italy'', ''Sweden'']'
$ws.Range("K235").Value = '[''Belgium'', ''Italy'', ''Sweden'']'
$ws.Range("K236").Value = '[''Belgium'', ''Republic of Ireland'', ''Italy'']'
$ws.Range("K237").Value = '[''Belgium'', ''Republic of Ireland'', ''Italy'']'
$ws.Range("K238").Value = '[''Italy'', ''Switzerland'', ''Wales'']'
$ws.Range("K239").Value = '[''Italy'', ''Switzerland'', ''Wales'']'
$ws.Range("K240").Value = '[''Italy'', ''Switzerland'', ''Wales'']'
$ws.Range("K241").Value = '[''Italy'', ''Switzerland'', ''Wales'']'
$ws.Range("K242").Value = '[''Italy'', ''Switzerland'', ''Wales'']'
$ws.Range("K243").Value = '[''Italy'', ''Switzerland'', ''Wales'']'
$ws.Range("K244").Value = '[''Italy'', ''Switzerland'', ''Wales'']'
$ws.Range("K245").Value = '[''Netherlands'', ''Austria'', ''Ukraine'']'
$ws.Range("K246").Value = '[''Netherlands'', ''Austria'', ''Ukraine'']'
$ws.Range("K247").Value = '[''Netherlands'', ''Austria'', ''Ukraine'']'
$ws.Range("K248").Value = '[''Netherlands'', ''Austria'', ''Ukraine'']'
$ws.Range("K249").Value = '[''Netherlands'', ''Austria'', ''Ukraine'']'
$ws.Range("K250").Value = '[''Netherlands'', ''Austria'', ''Ukraine'']'
$ws.Range("K251").Value = '[''Netherlands'', ''Austria'', ''Ukraine'']'
$ws.Range("K252").Value = '[''Belgium'', ''Finland'', ''Russia'']'
$ws.Range("K253").Value = '[''Belgium'', ''Finland'', ''Russia'']'
$ws.Range("K254").Value = '[''Belgium'', ''Denmark'', ''Russia'']'
$ws.Range("K255").Value = '[''Belgium'', ''Finland'', ''Denmark'']'
$ws.Range("K256").Value = '[''Belgium'', ''Finland'', ''Denmark'']'
$ws.Range("K257").Value = '[''Belgium'', ''Finland'', ''Denmark'']'
$ws.Range("K258").Value = '[''Belgium'', ''Finland'', ''Denmark'']'
$ws.Range("K259").Value = '[''Belgium'', ''Finland'', ''Denmark'']'
$ws.Range("K260").Value = '[''Belgium'', ''Finland'', ''Denmark'']'
$ws.Range("K261").Value = '[''England'', ''Croatia'', ''Czech Republic'']'
$ws.Range("K262").Value = '[''England'', ''Croatia'', ''Czech Republic'']'
$ws.Range("K263").Value = '[''England'', ''Croatia'', ''Czech Republic'']'
$ws.Range("K264").Value = '[''England'', ''Croatia'', ''Czech Republic'']'
$ws.Range("K265").Value = '[''England'', ''Croatia'', ''Czech Republic'']'
$ws.Range("K266").Value = '[''England'', ''Croatia'', ''Czech Republic'']'
$ws.Range("K267").Value = '[''England'', ''Croatia'', ''Czech Republic'']'
$ws.Range("K280").Value = '[''France'', ''Germany'', ''Portugal'']'
$ws.Range("K281").Value = '[''France'', ''Hungary'', ''Portugal'']'
$ws.Range("K282").Value = '[''France'', ''Hungary'', ''Portugal'']'
$ws.Range("K283").Value = '[''France'', ''Hungary'', ''Portugal'']'
$ws.Range("K285").Value = '[''France'', ''Hungary'', ''Portugal'']'
$ws.Range("K286").Value = '[''France'', ''Germany'', ''Portugal'']'
$ws.Range("K287").Value = '[''France'', ''Hungary'', ''Portugal'']'
$ws.Range("K288").Value = '[''France'', ''Germany'', ''Portugal'']'
$ws.Range("K289").Value = '[''Scotland'', ''Germany'', ''Switzerland'']'
$ws.Range("K290").Value = '[''Scotland'', ''Germany'', ''Switzerland'']'
$ws.Range("K291").Value = '[''Scotland'', ''Germany'', ''Switzerland'']'
$ws.Range("K292").Value = '[''Scotland'', ''Germany'', ''Switzerland'']'
$ws.Range("K293").Value = '[''Hungary'', ''Germany'', ''Switzerland'']'
$ws.Range("K294").Value = '[''Albania'', ''Spain'', ''Italy'']'
$ws.Range("K295").Value = '[''Albania'', ''Spain'', ''Italy'']'
$ws.Range("K296").Value = '[''Croatia'', ''Spain'', ''Italy'']'
$ws.Range("K297").Value = '[''Croatia'', ''Spain'', ''Italy'']'
$ws.Range("K298").Value = '[''Croatia'', ''Spain'', ''Italy'']'
$ws.Range("K307").Value = '[''England'', ''Slovenia'', ''Denmark'']'
$ws.Range("K308").Value = '[''England'', ''Slovenia'', ''Denmark'']'
$ws.Range("K309").Value = '[''England'', ''Slovenia'', ''Denmark'']'
$ws.Range("K310").Value = '[''Belgium'', ''Romania'', ''Ukraine'']'
$ws.Range("K311").Value = '[''Belgium'', ''Romania'', ''Slovakia'']'
$ws.Range("K312").Value = '[''Belgium'', ''Ukraine'', ''Slovakia'']'
$ws.Range("K313").Value = '[''Belgium'', ''Romania'', ''Slovakia'']'
$ws.Range("K314").Value = '[''Czech Republic'', ''Portugal'', ''Turkey'']'
$ws.Range("K315").Value = '[''Czech Republic'', ''Portugal'', ''Turkey'']'
$ws.Range("K316").Value = '[''Georgia'', ''Portugal'', ''Turkey'']'
$ws.Range("K317").Value = '[''Georgia'', ''Portugal'', ''Turkey'']'
$ws.Range("K318").Value = '[''Georgia'', ''Portugal'', ''Turkey'']'
$ws.Range("K319").Value = '[''Georgia'', ''Portugal'', ''Turkey'']'
$ws.Range("K320").Value = '[''Georgia'', ''Portugal'', ''Turkey'']'
